$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds the last-changed date serial for every data
# row (2 through 357). The update bumps that date by one day (45202 -> 45203,
# i.e. 2023-10-03 -> 2023-10-04) across the whole column.
$ws.Range("C2:C357").Value = 45203
